$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.568.97"
$ws.Range("E2").Value = "  +0.67%  "

$ws.Range("D3").Value = "2.695.76"
$ws.Range("E3").Value = "  +2.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "598.86"
$ws.Range("E5").Value = "  +0.43%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "160.18"
$ws.Range("E6").Value = "  +3.00%  "

$ws.Range("E7").Value = "  -0.03%  "

$ws.Range("E8").Value = "  +0.42%  "

$ws.Range("D9").Value = "2.694.92"
$ws.Range("E9").Value = "  +2.06%  "

$ws.Range("E10").Value = "  +0.19%  "

$ws.Range("E11").Value = "  -0.45%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.30"
$ws.Range("E12").Value = "  +1.31%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.360"
$ws.Range("E13").Value = "  +2.58%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.26"
$ws.Range("E14").Value = "  +1.27%  "

$ws.Range("D15").Value = "3.185.78"
$ws.Range("E15").Value = "  +2.02%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000189"
$ws.Range("E16").Value = "  -0.06%  "

$ws.Range("D17").Value = "68.407.35"
$ws.Range("E17").Value = "  +0.60%  "

$ws.Range("D18").Value = "2.694.33"
$ws.Range("E18").Value = "  +1.77%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.86"
$ws.Range("E19").Value = "  +4.80%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "367.54"
$ws.Range("E20").Value = "  +1.26%  "

$ws.Range("E21").Value = "  +3.40%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.53"
$ws.Range("E22").Value = "  +2.88%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.88"
$ws.Range("E23").Value = "  +2.45%  "

$ws.Range("E24").Value = "  +3.19%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "74.57"
$ws.Range("E25").Value = "  -0.01%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  -0.11%  "

$ws.Range("E27").Value = "  +3.10%  "

$ws.Range("D28").Value = "2.826.05"
$ws.Range("E28").Value = "  +1.88%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0000105"
$ws.Range("E29").Value = "  +0.89%  "

$ws.Range("E30").Value = "  +0.35%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "573.03"
$ws.Range("E31").Value = "  +3.44%  "

$ws.Range("E32").Value = "  +2.78%  "

$ws.Range("E33").Value = "  +3.53%  "

$ws.Range("E35").Value = "  +3.11%  "

$ws.Range("E36").Value = "  +6.62%  "

$ws.Range("E37").Value = "  -0.02%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "161.40"
$ws.Range("E38").Value = "  +0.17%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "19.88"
$ws.Range("E39").Value = "  +2.42%  "

$ws.Range("E40").Value = "  +2.09%  "

$ws.Range("E41").Value = "  +2.24%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.40"
$ws.Range("E42").Value = "  +1.87%  "

$ws.Range("E43").Value = "  +2.31%  "

$ws.Range("E44").Value = "  +0.34%  "

$ws.Range("D46").Value = "0.0₆0316"
$ws.Range("E46").Value = "  -5.64%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "157.49"
$ws.Range("E47").Value = "  -1.19%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.97"
$ws.Range("E48").Value = "  +6.94%  "

$ws.Range("E49").Value = "  +5.14%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.599"
$ws.Range("E50").Value = "  +6.80%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "22.02"
$ws.Range("E51").Value = "  +0.18%  "
